# Update the judge dashboard roster: two defendants' names were corrected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (21CRB01437-A): MANCHESTER, LINDSEY -> KUDELA, JUSTIN
$ws.Range("C10").Value = "KUDELA"
$ws.Range("D10").Value = "JUSTIN"

# Row 13 (21CRB01387-A): MURPHY, KEVIN -> NETTLER, KATHERINE
$ws.Range("C13").Value = "NETTLER"
$ws.Range("D13").Value = "KATHERINE"

# Move the view/selection like the author left it (scrolled down, cell D13 selected)
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D13").Select()
